# Weekly data refresh: a new price record (most recent week) is inserted
# as a new row 31 in the "Alcachofa" sheet, pushing all the existing
# records from row 31 downward by one row (old row 31 -> new row 32,
# old row 53 -> new row 54, etc.). The sheet dimension grows from
# A1:R53 to A1:R54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 31; this shifts rows 31:53 down to 32:54
# and keeps all previously existing data/formatting intact.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new weekly record.
$ws.Cells.Item(31, 1).Value = 11
$ws.Cells.Item(31, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(31, 3).Value = "Bíobío"
$ws.Cells.Item(31, 4).Value = 44762
$ws.Cells.Item(31, 5).Value = 8
$ws.Cells.Item(31, 6).Value = 100112013
$ws.Cells.Item(31, 7).Value = "Alcachofa"
$ws.Cells.Item(31, 8).Value = "Española"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 70
$ws.Cells.Item(31, 11).Value = 15000
$ws.Cells.Item(31, 12).Value = 16000
$ws.Cells.Item(31, 13).Value = 15571
$ws.Cells.Item(31, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(31, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value = 519
$ws.Cells.Item(31, 17).Value = 30
$ws.Cells.Item(31, 18).Value = "Hortaliza"
